$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Forecast Comparison": shift Week_Start_Date forward one week and
# replace the MyForecast values with the new "penalty/reward" figures.
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$newForecast = @(4, 4, 5, 5, 6, 4, 6, 5, 4, 4, 4, 4, 4, 3, 3, 3)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2

    # Column B: Week_Start_Date — force text so it stays an inline string
    # rather than being auto-parsed into a real Excel date serial number.
    $cellB = $wsForecast.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $newDates[$i]

    # Column D: MyForecast — plain numeric value.
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Summary": update the derived metrics that rolled with the new week.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# B2: Historical Range
$wsSummary.Range("B2").Value = "2023-01-29 to 2025-01-05"

# B8: Total Historical Sales
$wsSummary.Range("B8").Value = "3723 units"

# B9-B11: rolling forecast totals — keep as text, not auto-converted numbers.
$cellB9 = $wsSummary.Range("B9")
$cellB9.NumberFormat = "@"
$cellB9.Value = "68"

$cellB10 = $wsSummary.Range("B10")
$cellB10.NumberFormat = "@"
$cellB10.Value = "37"

$cellB11 = $wsSummary.Range("B11")
$cellB11.NumberFormat = "@"
$cellB11.Value = "17"

# B12: Max Forecast
$cellB12 = $wsSummary.Range("B12")
$cellB12.NumberFormat = "@"
$cellB12.Value = "6"

# B14: Min Forecast
$cellB14 = $wsSummary.Range("B14")
$cellB14.NumberFormat = "@"
$cellB14.Value = "3"

# B15: Min Forecast Week — text date, force text format.
$cellB15 = $wsSummary.Range("B15")
$cellB15.NumberFormat = "@"
$cellB15.Value = "2025-04-20"
